# Updating filtered feeds from workflow
# Appends three new feed rows (83-85) to the "Filtered Feeds" sheet,
# each with a link (A, hyperlinked), keywords (B) and title (C), matching
# the upstream workflow's latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell whose formatting (the built-in "Hyperlink" cell style used
# throughout column A) we reuse for the new link cells, so the new rows line
# up with the rest of the sheet instead of minting fresh style records.
$styleSource = $ws.Range("A2")

$newRows = @(
    @{
        Row   = 83
        Link  = "https://www.genomeweb.com/cancer/guardant-health-liquid-biopsy-nabs-fda-approval-cdx-pfizer-colorectal-cancer-drug-combo"
        Kw    = "CDx"
        Title = "Guardant Health Liquid Biopsy Nabs FDA Approval as CDx for Pfizer Colorectal Cancer Drug Combo"
    },
    @{
        Row   = 84
        Link  = "https://www.360dx.com/cancer/circulating-tumor-cell-assay-finds-best-responders-amgens-small-cell-lung-cancer-drug"
        Kw    = "companion diagnostic"
        Title = "Circulating Tumor Cell Assay Finds Best Responders to Amgen's Small Cell Lung Cancer Drug"
    },
    @{
        Row   = 85
        Link  = "https://www.360dx.com/cancer/guardant-health-liquid-biopsy-nabs-fda-approval-cdx-pfizer-colorectal-cancer-drug-combo"
        Kw    = "CDx"
        Title = "Guardant Health Liquid Biopsy Nabs FDA Approval as CDx for Pfizer Colorectal Cancer Drug Combo"
    }
)

foreach ($item in $newRows) {
    $r = $item.Row

    $linkCell = $ws.Cells.Item($r, 1)
    $kwCell   = $ws.Cells.Item($r, 2)
    $titleCell = $ws.Cells.Item($r, 3)

    # Column A: add the real hyperlink (sets the relationship + cell text).
    $ws.Hyperlinks.Add($linkCell, $item.Link)

    # Re-apply the shared "Hyperlink" look from an existing link cell so we
    # don't leave the freshly-added hyperlink using a brand-new style record.
    $styleSource.Copy()
    $linkCell.PasteSpecial(-4122)

    # Column B / C: plain text values.
    $kwCell.Value = $item.Kw
    $titleCell.Value = $item.Title
}

$excel.CutCopyMode = 0

Write-Output "Added rows 83-85 to Filtered Feeds"
